# "Generate Report for Handoff"
#
# The localization CI re-generated the handoff report: a new working file
# (f5e311f0-4e78-4d0a-acd0-c09944a03258.md) replaces the previous one
# (c0a39f1d-4e12-46f2-aa76-18eba71c7bed.md), new xliff package hashes were
# produced, and the handoff timestamps were refreshed.

$wb = $excel.ActiveWorkbook

# The hyperlink target (stored in the worksheet's relationship part) is left
# untouched by the report generator - only the on-sheet display text next to
# each link is refreshed to the new file name.
$linkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/df99ea405929caf0d925ee8e1a1176846275544c/e2e/c0a39f1d-4e12-46f2-aa76-18eba71c7bed.md"

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Range("A2").Value = "f5e311f0-4e78-4d0a-acd0-c09944a03258.md"
$wsOverview.Range("B2").Value = "e2e\f5e311f0-4e78-4d0a-acd0-c09944a03258.md"
$wsOverview.Range("G2").Value = "2017-01-03 07:54:56"

# --- Sheet "zh-cn" ---
$wsZhCn = $wb.Worksheets.Item(2)
$wsZhCn.Range("A2").Value = "f5e311f0-4e78-4d0a-acd0-c09944a03258.md"
$wsZhCn.Range("G2").Value = "f5e311f0-4e78-4d0a-acd0-c09944a03258.5238f77fb70b8c273bdb030176ed226884c7a133.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2017-01-03 07:54:45"

# --- Sheet "de-de" ---
$wsDeDe = $wb.Worksheets.Item(3)
$wsDeDe.Range("A2").Value = "f5e311f0-4e78-4d0a-acd0-c09944a03258.md"
$wsDeDe.Range("G2").Value = "f5e311f0-4e78-4d0a-acd0-c09944a03258.5238f77fb70b8c273bdb030176ed226884c7a133.de-de.xlf"
$wsDeDe.Range("H2").Value = "2017-01-03 07:54:56"

# --- Refresh the hyperlink display text on each sheet ---
# (Hyperlink property setters on this platform append a duplicate hyperlink
# rather than editing in place, so we clear each sheet's hyperlinks and
# re-add a single one pointing at the original, unchanged target URL.)
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $linkTarget, "", "", "e2e\f5e311f0-4e78-4d0a-acd0-c09944a03258.md") | Out-Null

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $linkTarget, "", "", "f5e311f0-4e78-4d0a-acd0-c09944a03258.md") | Out-Null

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $linkTarget, "", "", "f5e311f0-4e78-4d0a-acd0-c09944a03258.md") | Out-Null
